$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# Update sheet view top-left cell (scroll so row 20 is at the top)
$excel.ActiveWindow.ScrollRow = 20

# Row 22: set date, Changes text, Notes text
$ws.Range("A22").Value = "10/28/2025"
$ws.Range("B22").Value = "Changes`n- MODIFIED: MyMIF.mif, Deliverable4.mif, TestCase7.mif`n- COMPLETED: Finished testing all of the branch instructions using MyMIF.mif. Then tested with Deliverable4.mif and TestCase7.mif. Changed some things in those .mif files because of my design change to addressing. So everything works but some of the instructions work differently because of how I use addresses, but I could definitely change that so it works as intended. And the things with inport and outport don't work because it isn't running on an actual board                                                                                                                                                                                                                                      "
$ws.Range("D22").Value = "Notes`n- I did further change TestCase7.mif so it fits with my addressing so that it works as expected. The only part that doesn't work now is writing/reading to and from the outport/inport(s), but that was to be expected, as I wasn't planning on implementing that`n- I guess I'm done now lol. YAYYYYY!!!`nBugs`n- "

# Set row height to match wrapped-text autofit
$ws.Rows.Item(22).RowHeight = 129.6
